$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Tests: BinomioDeNewton y MiMath" increment row (row 24)
$ws.Range("C24").Value = "Tests: BinomioDeNewton y MiMath"
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 15/1440
$ws.Range("H24").Value = 22/24 + 37/1440
$ws.Range("I24").Value = 22/24 + 51/1440
$ws.Range("J24").Formula = '=IFERROR(IF(OR(ISBLANK(H24),ISBLANK(I24)),"",IF(I24>=H24,I24-H24,"Error")),"Error")'
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 90
$ws.Range("N24").Formula = '=IFERROR(IF(OR(J24="",ISBLANK(L24)),"",J24+L24),"Error")'

$ws.Range("M24").Select()
